$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds line1..line6 (rows 2-7) followed directly by
# extr1..extr8 (rows 8-15). Two new rows - "line7" and "line8" - are being
# inserted right after "line6", so the extr1..extr8 block needs to shift
# down two rows (old rows 8-15 -> new rows 10-17) before the new line7/
# line8 values are written into the now-vacated rows 8 and 9.

# Shift the extr1..extr8 block (8 rows) down by two rows, values+formats.
$ws.Range("A8:E15").Copy()
$ws.Range("A10").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Rows 16/17 did not exist before the shift, so the paste above did not
# carry over column A's "index" style onto them - copy it over explicitly
# from the row above (which already has it) so no new style gets created.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Overwrite rows 8 and 9 (still holding stale copies of the old extr1/
# extr2 rows) with the new line7 / line8 data.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber the "index" column for the extr1..extr8 rows, which shifted
# down by two rows (old row 8 -> new row 10, ... old row 15 -> new row 17).
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# Update the from_bus / to_bus values for the extr rows to their new
# targets (values shifted by two positions compared to the old sheet).
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7

$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
